$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.3

$ws.Range("B3").Value = 1.54
$ws.Range("D3").Value = 1.4
$ws.Range("E3").Value = 1.32
$ws.Range("F3").Value = 1.21

$ws.Range("C4").Value = 1.44
$ws.Range("E4").Value = 1.23

$ws.Range("C5").Value = 1.35
$ws.Range("D5").Value = 1.33

$ws.Range("C6").Value = 1.45

$ws.Range("G7").Value = 1.15
